$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "62.869.02"
$ws.Range("E2").Value = "  +0.04%  "
$ws.Range("D3").Value = "3.451.14"
$ws.Range("E3").Value = "  -0.60%  "
$ws.Range("E4").Value = "  -0.04%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "580.30"
$ws.Range("E5").Value = "  -0.03%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "150.50"
$ws.Range("E6").Value = "  +2.24%  "
$ws.Range("E7").Value = "  -0.01%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.488"
$ws.Range("E8").Value = "  +1.30%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "8.08"
$ws.Range("E9").Value = "  +6.23%  "
$ws.Range("E10").Value = "  -0.38%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.416"
$ws.Range("E11").Value = "  +3.97%  "
$ws.Range("D12").Value = "4.047.38"
$ws.Range("E12").Value = "  -0.47%  "
$ws.Range("E13").Value = "  -0.58%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "28.32"
$ws.Range("E14").Value = "  -4.48%  "
$ws.Range("D15").Value = "3.447.01"
$ws.Range("E15").Value = "  -0.41%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.0000173"
$ws.Range("D17").Value = "62.863.76"
$ws.Range("E17").Value = "  -0.02%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "6.43"
$ws.Range("E18").Value = "  +0.97%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "14.60"
$ws.Range("E19").Value = "  +1.73%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "8.99"
$ws.Range("E20").Value = "  -2.59%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "388.47"
$ws.Range("E21").Value = "  +0.01%  "
$ws.Range("E22").Value = "  +1.02%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "75.27"
$ws.Range("E23").Value = "  +0.73%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "1.00"
$ws.Range("E24").Value = "  +0.01%  "
$ws.Range("E25").Value = "  +0.50%  "
$ws.Range("D26").Value = "3.587.71"
$ws.Range("E26").Value = "  -0.57%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "0.185"
$ws.Range("E27").Value = "  +3.15%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "7.75"
$ws.Range("E28").Value = "  +1.92%  "
$ws.Range("E29").Value = "  +0.13%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "8.03"
$ws.Range("E30").Value = "  -1.67%  "
$ws.Range("E31").Value = "  -0.66%  "
$ws.Range("E32").Value = "  -0.03%  "
$ws.Range("E33").Value = "  -2.51%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "23.33"
$ws.Range("E34").Value = "  -1.76%  "
$ws.Range("E35").Value = "  +2.86%  "
$ws.Range("E36").Value = "  +2.97%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "32.13"
$ws.Range("E37").Value = "  +2.22%  "
$ws.Range("E38").Value = "  -1.65%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "168.98"
$ws.Range("E39").Value = "  -0.83%  "
$ws.Range("D40").Value = "3.491.01"
$ws.Range("E40").Value = "  -0.56%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.0782"
$ws.Range("E41").Value = "  +1.88%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "42.85"
$ws.Range("E42").Value = "  +1.45%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.786"
$ws.Range("E43").Value = "  -1.82%  "
$ws.Range("E44").Value = "  -0.88%  "
$ws.Range("E45").Value = "  -2.04%  "
$ws.Range("E46").Value = "  -1.25%  "
$ws.Range("D47").Value = "2.559.74"
$ws.Range("E47").Value = "  -1.94%  "
$ws.Range("E48").Value = "  +2.58%  "
$ws.Range("E49").Value = "  +0.87%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "22.86"
$ws.Range("E50").Value = "  -2.14%  "
$ws.Range("E51").Value = "  +0.05%  "
